$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.6826774943282885
$ws.Cells.Item(2, 4).Value = 0.09158574049555668
$ws.Cells.Item(2, 5).Value = 0.1253963583736044
$ws.Cells.Item(2, 6).Value = 2.234179487570046
$ws.Cells.Item(2, 7).Value = 1.572035072809442
$ws.Cells.Item(2, 8).Value = 1.417889399802348
$ws.Cells.Item(2, 9).Value = 1.483147664041056
$ws.Cells.Item(2, 10).Value = 0.1857805617549593
$ws.Cells.Item(2, 11).Value = 0.5960625391665531
$ws.Cells.Item(2, 12).Value = 0.3146078491034245
$ws.Cells.Item(3, 2).Value = 0.660473036343916
$ws.Cells.Item(3, 4).Value = 0.09138973107846482
$ws.Cells.Item(3, 5).Value = 0.1257938579253786
$ws.Cells.Item(3, 6).Value = 2.237321918419767
$ws.Cells.Item(3, 7).Value = 1.571585008173258
$ws.Cells.Item(3, 8).Value = 1.423074462990186
$ws.Cells.Item(3, 9).Value = 1.493325324341253
$ws.Cells.Item(3, 10).Value = 0.1864827754605707
$ws.Cells.Item(3, 11).Value = 0.5407404951086505
$ws.Cells.Item(3, 12).Value = 0.3040251252710533
$ws.Cells.Item(4, 2).Value = 0.6471118241173599
$ws.Cells.Item(4, 4).Value = 0.09128477558207138
$ws.Cells.Item(4, 5).Value = 0.1260566787596358
$ws.Cells.Item(4, 6).Value = 2.240303357876243
$ws.Cells.Item(4, 7).Value = 1.572170734290765
$ws.Cells.Item(4, 8).Value = 1.426851160199718
$ws.Cells.Item(4, 9).Value = 1.500151728108172
$ws.Cells.Item(4, 10).Value = 0.1869388829846006
$ws.Cells.Item(4, 11).Value = 0.5069297553352783
$ws.Cells.Item(4, 12).Value = 0.2976733564793079
$ws.Cells.Item(5, 2).Value = 0.6417360114998303
$ws.Cells.Item(5, 4).Value = 0.09124589731548483
$ws.Cells.Item(5, 5).Value = 0.1261685081208199
$ws.Cells.Item(5, 6).Value = 2.241782893658787
$ws.Cells.Item(5, 7).Value = 1.572626046879037
$ws.Cells.Item(5, 8).Value = 1.42853937672993
$ws.Cells.Item(5, 9).Value = 1.503078739061593
$ws.Cells.Item(5, 10).Value = 0.1871310368908445
$ws.Cells.Item(5, 11).Value = 0.4931916191893038
$ws.Cells.Item(5, 12).Value = 0.2951218111151661
$ws.Cells.Item(6, 2).Value = 0.6408475431730949
$ws.Cells.Item(6, 4).Value = 0.09123967728751126
$ws.Cells.Item(6, 5).Value = 0.1261873631916335
$ws.Cells.Item(6, 6).Value = 2.242044550565666
$ws.Cells.Item(6, 7).Value = 1.572714730422675
$ws.Cells.Item(6, 8).Value = 1.428828715653765
$ws.Cells.Item(6, 9).Value = 1.503573537501662
$ws.Cells.Item(6, 10).Value = 0.1871633239844546
$ws.Cells.Item(6, 11).Value = 0.4909128432373677
$ws.Cells.Item(6, 12).Value = 0.2947003583538645
$ws.Cells.Item(7, 2).Value = 0.6470390440406959
$ws.Cells.Item(7, 4).Value = 0.09128423546953712
$ws.Cells.Item(7, 5).Value = 0.1260581677709847
$ws.Cells.Item(7, 6).Value = 2.240322240100994
$ws.Cells.Item(7, 7).Value = 1.572175997893879
$ws.Cells.Item(7, 8).Value = 1.426873323978398
$ws.Cells.Item(7, 9).Value = 1.500190614866067
$ws.Cells.Item(7, 10).Value = 0.1869414489677146
$ws.Cells.Item(7, 11).Value = 0.5067443153664328
$ws.Cells.Item(7, 12).Value = 0.2976387960704301
$ws.Cells.Item(8, 2).Value = 0.6749651416974132
$ws.Cells.Item(8, 4).Value = 0.09151497341321857
$ws.Cells.Item(8, 5).Value = 0.1255295307754025
$ws.Cells.Item(8, 6).Value = 2.235044684244045
$ws.Cells.Item(8, 7).Value = 1.571700874042776
$ws.Cells.Item(8, 8).Value = 1.419554153879446
$ws.Cells.Item(8, 9).Value = 1.486537127620256
$ws.Cells.Item(8, 10).Value = 0.1860175158460984
$ws.Cells.Item(8, 11).Value = 0.5769551369615158
$ws.Cells.Item(8, 12).Value = 0.3109286789628811
$ws.Cells.Item(9, 2).Value = 0.7318715974087979
$ws.Cells.Item(9, 4).Value = 0.09208876202869476
$ws.Cells.Item(9, 5).Value = 0.1246411640642169
$ws.Cells.Item(9, 6).Value = 2.233042500096118
$ws.Cells.Item(9, 7).Value = 1.577618097182224
$ws.Cells.Item(9, 8).Value = 1.409905272998103
$ws.Cells.Item(9, 9).Value = 1.464341783530983
$ws.Cells.Item(9, 10).Value = 0.1844029917320003
$ws.Cells.Item(9, 11).Value = 0.7158720901717004
$ws.Cells.Item(9, 12).Value = 0.338145395119426
$ws.Cells.Item(10, 2).Value = 0.7749683489082884
$ws.Cells.Item(10, 4).Value = 0.09258325248898558
$ws.Cells.Item(10, 5).Value = 0.1240781788812662
$ws.Cells.Item(10, 6).Value = 2.236663615438204
$ws.Cells.Item(10, 7).Value = 1.586155990619929
$ws.Cells.Item(10, 8).Value = 1.405683061254621
$ws.Cells.Item(10, 9).Value = 1.450824846720771
$ws.Cells.Item(10, 10).Value = 0.1833362052811658
$ws.Cells.Item(10, 11).Value = 0.8186790119964655
$ws.Cells.Item(10, 12).Value = 0.358843392535519
$ws.Cells.Item(11, 2).Value = 0.7948500636063898
$ws.Cells.Item(11, 4).Value = 0.09282383761232893
$ws.Cells.Item(11, 5).Value = 0.123841391742566
$ws.Cells.Item(11, 6).Value = 2.239417661314434
$ws.Cells.Item(11, 7).Value = 1.590953506540416
$ws.Cells.Item(11, 8).Value = 1.404384703262735
$ws.Cells.Item(11, 9).Value = 1.445281172546245
$ws.Cells.Item(11, 10).Value = 0.1828766374303228
$ws.Cells.Item(11, 11).Value = 0.8656093924886648
$ws.Cells.Item(11, 12).Value = 0.3684115142327187
$ws.Cells.Item(12, 2).Value = 0.8024181390495073
$ws.Cells.Item(12, 4).Value = 0.09291717125560339
$ws.Cells.Item(12, 5).Value = 0.1237544928762064
$ws.Cells.Item(12, 6).Value = 2.240619727145798
$ws.Cells.Item(12, 7).Value = 1.592901800530399
$ws.Cells.Item(12, 8).Value = 1.403982519841691
$ws.Cells.Item(12, 9).Value = 1.443268943143103
$ws.Cells.Item(12, 10).Value = 0.1827062954606182
$ws.Cells.Item(12, 11).Value = 0.8834038343643726
$ws.Cells.Item(12, 12).Value = 0.3720565604266994
$ws.Cells.Item(13, 2).Value = 0.8007864785726326
$ws.Cells.Item(13, 4).Value = 0.0928969713471588
$ws.Cells.Item(13, 5).Value = 0.1237730851936476
$ws.Cells.Item(13, 6).Value = 2.240353761464846
$ws.Cells.Item(13, 7).Value = 1.592476345268977
$ws.Cells.Item(13, 8).Value = 1.404065157990175
$ws.Cells.Item(13, 9).Value = 1.443698441538267
$ws.Cells.Item(13, 10).Value = 0.1827428178863191
$ws.Cells.Item(13, 11).Value = 0.8795704714400472
$ws.Cells.Item(13, 12).Value = 0.3712705666231386
$ws.Cells.Item(14, 2).Value = 0.79547190891887
$ws.Cells.Item(14, 4).Value = 0.09283147166019035
$ws.Cells.Item(14, 5).Value = 0.1238341871201465
$ws.Cells.Item(14, 6).Value = 2.239513365936844
$ws.Cells.Item(14, 7).Value = 1.591111155641727
$ws.Cells.Item(14, 8).Value = 1.404349822236668
$ws.Cells.Item(14, 9).Value = 1.44511388078304
$ws.Cells.Item(14, 10).Value = 0.1828625494806029
$ws.Cells.Item(14, 11).Value = 0.8670728953984792
$ws.Cells.Item(14, 12).Value = 0.3687109580493342
$ws.Cells.Item(15, 2).Value = 0.7922216849898973
$ws.Cells.Item(15, 4).Value = 0.09279164095276116
$ws.Cells.Item(15, 5).Value = 0.1238719738847216
$ws.Cells.Item(15, 6).Value = 2.23901932827637
$ws.Cells.Item(15, 7).Value = 1.590292079148384
$ws.Cells.Item(15, 8).Value = 1.404535839332084
$ws.Cells.Item(15, 9).Value = 1.445992213504567
$ws.Cells.Item(15, 10).Value = 0.1829363682657583
$ws.Cells.Item(15, 11).Value = 0.8594207396061222
$ws.Cells.Item(15, 12).Value = 0.3671459600642493
$ws.Cells.Item(16, 2).Value = 0.7736745589358804
$ws.Cells.Item(16, 4).Value = 0.09256784270022322
$ws.Cells.Item(16, 5).Value = 0.124094041292298
$ws.Cells.Item(16, 6).Value = 2.236505904577797
$ws.Cells.Item(16, 7).Value = 1.585860862985399
$ws.Cells.Item(16, 8).Value = 1.405780433815437
$ws.Cells.Item(16, 9).Value = 1.451199318487625
$ws.Cells.Item(16, 10).Value = 0.1833667556251517
$ws.Cells.Item(16, 11).Value = 0.8156152371758481
$ws.Cells.Item(16, 12).Value = 0.3582211514109588
$ws.Cells.Item(17, 2).Value = 0.7623670327594141
$ws.Cells.Item(17, 4).Value = 0.09243454120316841
$ws.Cells.Item(17, 5).Value = 0.1242352126144901
$ws.Cells.Item(17, 6).Value = 2.235247484189458
$ws.Cells.Item(17, 7).Value = 1.583376593366737
$ws.Cells.Item(17, 8).Value = 1.406703338483155
$ws.Cells.Item(17, 9).Value = 1.454548731276937
$ws.Cells.Item(17, 10).Value = 0.1836373629854042
$ws.Cells.Item(17, 11).Value = 0.7887833608361632
$ws.Cells.Item(17, 12).Value = 0.3527850411628322
$ws.Cells.Item(18, 2).Value = 0.7558893249049561
$ws.Cells.Item(18, 4).Value = 0.09235934287135805
$ws.Cells.Item(18, 5).Value = 0.1243182295203429
$ws.Cells.Item(18, 6).Value = 2.234627848740942
$ws.Cells.Item(18, 7).Value = 1.582033685818246
$ws.Cells.Item(18, 8).Value = 1.407292748008743
$ws.Cells.Item(18, 9).Value = 1.456532193367423
$ws.Cells.Item(18, 10).Value = 0.1837954306487886
$ws.Cells.Item(18, 11).Value = 0.7733657476124165
$ws.Cells.Item(18, 12).Value = 0.3496726955318508
$ws.Cells.Item(19, 2).Value = 0.7537005792418938
$ws.Cells.Item(19, 4).Value = 0.09233413562718695
$ws.Cells.Item(19, 5).Value = 0.1243466503385116
$ws.Cells.Item(19, 6).Value = 2.234435942501392
$ws.Cells.Item(19, 7).Value = 1.581593761027193
$ws.Cells.Item(19, 8).Value = 1.407502373490757
$ws.Cells.Item(19, 9).Value = 1.457213543771992
$ws.Cells.Item(19, 10).Value = 0.1838493658977927
$ws.Cells.Item(19, 11).Value = 0.7681482646196969
$ws.Cells.Item(19, 12).Value = 0.3486213790929469
$ws.Cells.Item(20, 2).Value = 0.7635680424534996
$ws.Cells.Item(20, 4).Value = 0.09244857905144599
$ws.Cells.Item(20, 5).Value = 0.1242199965145794
$ws.Cells.Item(20, 6).Value = 2.235370663391777
$ws.Cells.Item(20, 7).Value = 1.583632148144119
$ws.Cells.Item(20, 8).Value = 1.406599031097102
$ws.Cells.Item(20, 9).Value = 1.45418628450377
$ws.Cells.Item(20, 10).Value = 0.1836083058439797
$ws.Cells.Item(20, 11).Value = 0.7916380726298939
$ws.Cells.Item(20, 12).Value = 0.353362238943248
$ws.Cells.Item(21, 2).Value = 0.7970318641898189
$ws.Cells.Item(21, 4).Value = 0.0928506501583044
$ws.Cells.Item(21, 5).Value = 0.1238161649837037
$ws.Cells.Item(21, 6).Value = 2.239755890571018
$ws.Cells.Item(21, 7).Value = 1.591508572271579
$ws.Cells.Item(21, 8).Value = 1.40426378120344
$ws.Cells.Item(21, 9).Value = 1.444695769685822
$ws.Cells.Item(21, 10).Value = 0.1828272814556076
$ws.Cells.Item(21, 11).Value = 0.8707431183099459
$ws.Cells.Item(21, 12).Value = 0.3694621862955074
$ws.Cells.Item(22, 2).Value = 0.8191312559575863
$ws.Cells.Item(22, 4).Value = 0.09312641193369942
$ws.Cells.Item(22, 5).Value = 0.1235683627468207
$ws.Cells.Item(22, 6).Value = 2.243549652590858
$ws.Cells.Item(22, 7).Value = 1.597423308682835
$ws.Cells.Item(22, 8).Value = 1.4032590914656
$ws.Cells.Item(22, 9).Value = 1.439000506531109
$ws.Cells.Item(22, 10).Value = 0.1823383185902896
$ws.Cells.Item(22, 11).Value = 0.9225762297484152
$ws.Cells.Item(22, 12).Value = 0.3801114792474181
$ws.Cells.Item(23, 2).Value = 0.8073156192964177
$ws.Cells.Item(23, 4).Value = 0.09297805094447043
$ws.Cells.Item(23, 5).Value = 0.123699147491048
$ws.Cells.Item(23, 6).Value = 2.241439955138461
$ws.Cells.Item(23, 7).Value = 1.594196250333738
$ws.Cells.Item(23, 8).Value = 1.403747598288987
$ws.Cells.Item(23, 9).Value = 1.441993752300661
$ws.Cells.Item(23, 10).Value = 0.182597325600538
$ws.Cells.Item(23, 11).Value = 0.8948998858202515
$ws.Cells.Item(23, 12).Value = 0.3744161667744095
$ws.Cells.Item(24, 2).Value = 0.7630249938565612
$ws.Cells.Item(24, 4).Value = 0.09244222805803659
$ws.Cells.Item(24, 5).Value = 0.1242268699330742
$ws.Cells.Item(24, 6).Value = 2.235314650590297
$ws.Cells.Item(24, 7).Value = 1.583516346034784
$ws.Cells.Item(24, 8).Value = 1.406646005238414
$ws.Cells.Item(24, 9).Value = 1.454349966516681
$ws.Cells.Item(24, 10).Value = 0.1836214348150838
$ws.Cells.Item(24, 11).Value = 0.7903474312589367
$ws.Cells.Item(24, 12).Value = 0.3531012474912387
$ws.Cells.Item(25, 2).Value = 0.7162494108791577
$ws.Cells.Item(25, 4).Value = 0.09192065697177298
$ws.Cells.Item(25, 5).Value = 0.1248656884857331
$ws.Cells.Item(25, 6).Value = 2.232690170561753
$ws.Cells.Item(25, 7).Value = 1.575282533544964
$ws.Cells.Item(25, 8).Value = 1.412012048700802
$ws.Cells.Item(25, 9).Value = 1.469856150917153
$ws.Cells.Item(25, 10).Value = 0.1848187315557701
$ws.Cells.Item(25, 11).Value = 0.6781600846897504
$ws.Cells.Item(25, 12).Value = 0.3306591244382844
